$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "G713IC-HX008" was removed from the tracked product codes list. The
# remaining codes shift up one row, leaving the last row (5) blank.
$ws.Range("A2").Value = "NH.Q8QEX.006"
$ws.Range("A3").Value = "FX506HEB-HN148"
$ws.Range("A4").Value = "DL341015898541WP "
$ws.Range("A5").ClearContents()

# Touch row 5's properties (no-op outline level) so the now-empty row is
# still written out as a row element instead of being dropped entirely.
$ws.Rows.Item(5).OutlineLevel = 0

$ws.Range("B10").Select()
